# Patch БД 2.2.1 - SBRFNDFL-3844
# Dedupe/rename shared strings (yo-fication: Расчет->Расчёт, etc.),
# re-home the title string, and update view selections / active sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Заголовок")
$ws2 = $wb.Worksheets.Item("1.Реквизиты")
$ws3 = $wb.Worksheets.Item("2. Свед о дох")
$ws4 = $wb.Worksheets.Item("3. Свед о вычет")
$ws5 = $wb.Worksheets.Item("4. Аванс платеж")

# Заголовок
$ws1.Range("G2").Value = "ф. РНУ НДФЛ"
$ws1.Range("B4").Value = "Регистр налогового учёта на доходы физических лиц"
$ws1.Range("B5").Value = "за"

# 1.Реквизиты
$ws2.Range("B2").Value = "1"
$ws2.Range("C2").Value = "2"
$ws2.Range("D2").Value = "3"
$ws2.Range("E2").Value = "4"
$ws2.Range("F2").Value = "5"
$ws2.Range("G2").Value = "6"
$ws2.Range("H2").Value = "7"
$ws2.Range("I2").Value = "8"
$ws2.Range("J2").Value = "9"
$ws2.Range("K2").Value = "10"
$ws2.Range("L2").Value = "11"
$ws2.Range("M2").Value = "12"
$ws2.Range("N2").Value = "13"
$ws2.Range("O2").Value = "14"
$ws2.Range("P2").Value = "15"
$ws2.Range("Q2").Value = "16"
$ws2.Range("R2").Value = "17"
$ws2.Range("S2").Value = "18"
$ws2.Range("T2").Value = "19"
$ws2.Range("U2").Value = "20"
$ws2.Range("V2").Value = "21"
$ws2.Range("B3").Value = "№ п/п"
$ws2.Range("C3").Value = "Налогоплательщик.`nИНП"
$ws2.Range("D3").Value = "Налогоплательщик.`nФамилия"
$ws2.Range("E3").Value = "Налогоплательщик.`nИмя"
$ws2.Range("F3").Value = "Налогоплательщик.`nОтчество"
$ws2.Range("G3").Value = "Налогоплательщик.`nДата рождения"
$ws2.Range("H3").Value = "Гражданство`n(код страны)"
$ws2.Range("I3").Value = "ИНН.В Российской`nФедерации"
$ws2.Range("J3").Value = "ИНН.В стране`nгражданства"
$ws2.Range("K3").Value = "Документ, удостов.`nличность.Код"
$ws2.Range("L3").Value = "Документ, удостов.`nличность.Номер"
$ws2.Range("M3").Value = "Статус`n(код)"
$ws2.Range("N3").Value = "Код субъекта"
$ws2.Range("O3").Value = "Индекс"
$ws2.Range("P3").Value = "Район"
$ws2.Range("Q3").Value = "Город"
$ws2.Range("R3").Value = "Населенный пункт"
$ws2.Range("S3").Value = "Улица"
$ws2.Range("T3").Value = "Дом"
$ws2.Range("U3").Value = "Корпус"
$ws2.Range("V3").Value = "Квартира"
$ws2.Range("W3").Value = "СНИЛС"

# 2. Свед о дох
$ws3.Range("B2").Value = "1"
$ws3.Range("C2").Value = "2"
$ws3.Range("D2").Value = "3"
$ws3.Range("E2").Value = "4"
$ws3.Range("F2").Value = "5"
$ws3.Range("G2").Value = "6"
$ws3.Range("H2").Value = "7"
$ws3.Range("I2").Value = "8"
$ws3.Range("J2").Value = "9"
$ws3.Range("K2").Value = "10"
$ws3.Range("L2").Value = "11"
$ws3.Range("M2").Value = "12"
$ws3.Range("N2").Value = "13"
$ws3.Range("O2").Value = "14"
$ws3.Range("P2").Value = "15"
$ws3.Range("Q2").Value = "16"
$ws3.Range("R2").Value = "17"
$ws3.Range("S2").Value = "18"
$ws3.Range("T2").Value = "19"
$ws3.Range("U2").Value = "20"
$ws3.Range("V2").Value = "21"
$ws3.Range("W2").Value = "22"
$ws3.Range("X2").Value = "23"
$ws3.Range("Y2").Value = "24"
$ws3.Range("B3").Value = "№ п/п"
$ws3.Range("C3").Value = "ИНП"
$ws3.Range("D3").Value = "ID`nоперации"
$ws3.Range("E3").Value = "Доход.Вид.`nКод"
$ws3.Range("F3").Value = "Доход.Вид.`nПризнак"
$ws3.Range("G3").Value = "Доход.Дата.`nНачисление"
$ws3.Range("H3").Value = "Доход.Дата.`nВыплата"
$ws3.Range("I3").Value = "Доход.Источник выплаты. КПП"
$ws3.Range("J3").Value = "Доход.Источник выплаты. ОКТМО"
$ws3.Range("K3").Value = "Доход.Сумма.`nНачисление"
$ws3.Range("L3").Value = "Доход.Сумма.`nВыплата"
$ws3.Range("M3").Value = "Сумма  вычета"
$ws3.Range("N3").Value = "Налоговая база"
$ws3.Range("O3").Value = "НДФЛ.`nПроцентная ставка`n(%)"
$ws3.Range("P3").Value = "НДФЛ.Расчёт.`nДата"
$ws3.Range("Q3").Value = "НДФЛ.Расчёт.`nСумма.`nИсчисленный"
$ws3.Range("R3").Value = "НДФЛ.Расчёт.`nСумма.`nУдержанный"
$ws3.Range("S3").Value = "НДФЛ.Расчёт.`nСумма.`nНе удержанный"
$ws3.Range("T3").Value = "НДФЛ.Расчёт.`nСумма.Излишне удержанный"
$ws3.Range("U3").Value = "НДФЛ.Расчёт.`nСумма.Возвращённый налогоплательщику"
$ws3.Range("V3").Value = "Перечисление в бюджет.`nСрок"
$ws3.Range("W3").Value = "Перечисление в бюджет.`nПлатёжное поручение.`nДата"
$ws3.Range("X3").Value = "Перечисление в бюджет.`nПлатёжное поручение.`nНомер"
$ws3.Range("Y3").Value = "Перечисление в бюджет.`nПлатёжное поручение.`nСумма"

# 3. Свед о вычет
$ws4.Range("B2").Value = "1"
$ws4.Range("C2").Value = "2"
$ws4.Range("D2").Value = "3"
$ws4.Range("E2").Value = "4"
$ws4.Range("F2").Value = "5"
$ws4.Range("G2").Value = "6"
$ws4.Range("H2").Value = "7"
$ws4.Range("I2").Value = "8"
$ws4.Range("J2").Value = "9"
$ws4.Range("K2").Value = "10"
$ws4.Range("L2").Value = "11"
$ws4.Range("M2").Value = "12"
$ws4.Range("N2").Value = "13"
$ws4.Range("O2").Value = "14"
$ws4.Range("P2").Value = "15"
$ws4.Range("Q2").Value = "16"
$ws4.Range("B3").Value = "№ п/п"
$ws4.Range("C3").Value = "ИНП"
$ws4.Range("D3").Value = "Код вычета"
$ws4.Range("E3").Value = "Документ о праве`nна налоговый вычет.`nТип"
$ws4.Range("F3").Value = "Документ о праве`nна налоговый вычет.`nДата"
$ws4.Range("G3").Value = "Документ о праве`nна налоговый вычет.`nНомер"
$ws4.Range("H3").Value = "Документ о праве`nна налоговый вычет.`nКод источника"
$ws4.Range("I3").Value = "Документ о праве`nна налоговый вычет.`nСумма"
$ws4.Range("J3").Value = "Начисленный доход.`nID операции"
$ws4.Range("K3").Value = "Начисленный доход.`nДата"
$ws4.Range("L3").Value = "Начисленный доход.`nКод дохода"
$ws4.Range("M3").Value = "Начисленный доход.`nСумма"
$ws4.Range("N3").Value = "Применение вычета.Дата заявления о применении налогового вычета"
$ws4.Range("O3").Value = "Применение вычета.Сумма применённого вычета с начала налогового периода"
$ws4.Range("P3").Value = "Применение вычета.`nТекущий период.`nДата"
$ws4.Range("Q3").Value = "Применение вычета.`nТекущий период.`nСумма"

# 4. Аванс платеж
$ws5.Range("B2").Value = "1"
$ws5.Range("C2").Value = "2"
$ws5.Range("D2").Value = "3"
$ws5.Range("E2").Value = "4"
$ws5.Range("F2").Value = "5"
$ws5.Range("G2").Value = "6"
$ws5.Range("B3").Value = "№ п/п"
$ws5.Range("C3").Value = "ИНП"
$ws5.Range("D3").Value = "ID операции"
$ws5.Range("E3").Value = "Сумма фиксированного`nавансового платежа"
$ws5.Range("F3").Value = "Уведомление, подтверждающее`nправо на уменьшение налога на`nфиксированные авансовые платежи.`nНомер"
$ws5.Range("G3").Value = "Уведомление, подтверждающее`nправо на уменьшение налога на`nфиксированные авансовые платежи.`nДата выдачи уведомления"
$ws5.Range("H3").Value = "Уведомление, подтверждающее`nправо на уменьшение налога на`nфиксированные авансовые платежи.`nКод налогового органа, выдавшего уведомление"

# Update selections on each sheet (order matters: the final
# Select()/Activate() call determines which sheet ends up active).
$ws2.Range("C3").Select()
$ws3.Range("D3").Select()
$ws4.Range("D3").Select()
$ws5.Range("B2").Select()
$ws1.Range("B2").Select()

# "Заголовок" is the active tab in the target workbook.
$ws1.Activate()
$ws1.Range("B2").Select()
